# Adds torque-as-function-of-power entries (English + SI units) to the
# AeroToolKitFunctionList sheet, refreshes the text of the pre-existing
# "power as function of torque" rows (366-377) whose descriptions had to be
# re-worded/renumbered, and highlights the POWER category header in yellow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Refresh text on the existing POWER rows (366-377)
# ---------------------------------------------------------------------
$ws.Range("I366").Value = "Power in lbf-ft/sec, f(torque in lbf-ft, RPM)"
$ws.Range("J366").Value = "lbf-ft, RPM"
$ws.Range("K366").Value = "lbf-ft/sec"
$ws.Range("I367").Value = "Power in lbf-ft/sec, f(torque in lbf-ft, RPS)"
$ws.Range("J367").Value = "lbf-ft, RPS"
$ws.Range("K367").Value = "lbf-ft/sec"
$ws.Range("I368").Value = "Power in lbf-ft/sec, f(torque in lbf-inch, RPM)"
$ws.Range("J368").Value = "lbf-inch, RPM"
$ws.Range("K368").Value = "lbf-ft/sec"
$ws.Range("I369").Value = "Power in lbf-ft/sec, f(torque in lbf-inch, RPS)"
$ws.Range("J369").Value = "lbf-inch, RPS"
$ws.Range("K369").Value = "lbf-ft/sec"
$ws.Range("I370").Value = "Power in horsepower, f(torque in lbf-ft, RPM)"
$ws.Range("J370").Value = "lbf-ft, RPM"
$ws.Range("K370").Value = "horsepower"
$ws.Range("I371").Value = "Power in horsepower, f(torque in lbf-ft, RPS)"
$ws.Range("J371").Value = "lbf-ft, RPS"
$ws.Range("K371").Value = "horsepower"
$ws.Range("I372").Value = "Power in horsepower, f(torque in lbf-inch, RPM)"
$ws.Range("J372").Value = "lbf-inch, RPM"
$ws.Range("K372").Value = "horsepower"
$ws.Range("I373").Value = "Power in horsepower, f(torque in lbf-inch, RPS)"
$ws.Range("J373").Value = "lbf-inch, RPS"
$ws.Range("K373").Value = "horsepower"
$ws.Range("D374").Value = "torqueNewtonMeter"
$ws.Range("I374").Value = "Power in watts, f(torque in newton-meter, RPM)"
$ws.Range("J374").Value = "newton-meter, RPM"
$ws.Range("K374").Value = "watts"
$ws.Range("D375").Value = "torqueNewtonMeter"
$ws.Range("I375").Value = "Power in watts, f(torque in newton-meter, RPS)"
$ws.Range("J375").Value = "newton-meter, RPS"
$ws.Range("K375").Value = "watts"
$ws.Range("D376").Value = "torqueNewtonMeter"
$ws.Range("I376").Value = "Power in kilowatts, f(torque in newton-meter, RPM)"
$ws.Range("J376").Value = "newton-meter, RPM"
$ws.Range("K376").Value = "kilowatts"
$ws.Range("D377").Value = "torqueNewtonMeter"
$ws.Range("I377").Value = "Power in kilowatts, f(torque in newton-meter, RPS)"
$ws.Range("J377").Value = "newton-meter, RPS"
$ws.Range("K377").Value = "kilowatts"

# ---------------------------------------------------------------------
# 2) Highlight the POWER category header cells (column A, rows 366-377)
#    in yellow (was orange/FFC000) to flag the refreshed section.
# ---------------------------------------------------------------------
$ws.Range("A366:A377").Interior.Color = 65535

# ---------------------------------------------------------------------
# 3) Add the new TORQUE category rows (378-389): torque as f(power, rev)
#    Column formats (A/B category-header colors) are copied from an
#    existing category block so no new styles/fills get created.
# ---------------------------------------------------------------------
$ws.Range("A295").Copy() | Out-Null
$ws.Range("A378:A389").PasteSpecial(-4122) | Out-Null
$ws.Range("B84").Copy() | Out-Null
$ws.Range("B378:B389").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A378").Value = "TORQUE"
$ws.Range("B378").Value = "TORQUE"
$ws.Range("C378").Value = "AeroTorque_lbfFt_fPowerFtLbfPerSecRevPerMin"
$ws.Range("D378").Value = "powerFtLbfPerSec"
$ws.Range("E378").Value = " revPerMin"
$ws.Range("H378").Value = 2
$ws.Range("I378").Formula = '="Torque in " & K378 & ", f(" & J378 & ")"'
$ws.Range("J378").Value = "lbf-ft/sec, RPM"
$ws.Range("K378").Value = "lbf-ft"

$ws.Range("A379").Value = "TORQUE"
$ws.Range("B379").Value = "TORQUE"
$ws.Range("C379").Value = "AeroTorque_lbfFt_fPowerFtLbfPerSecRevPerSec"
$ws.Range("D379").Value = "powerFtLbfPerSec"
$ws.Range("E379").Value = " revPerSec"
$ws.Range("H379").Value = 2
$ws.Range("I379").Formula = '="Torque in " & K379 & ", f(" & J379 & ")"'
$ws.Range("J379").Value = "lbf-ft/sec, RPS"
$ws.Range("K379").Value = "lbf-ft"

$ws.Range("A380").Value = "TORQUE"
$ws.Range("B380").Value = "TORQUE"
$ws.Range("C380").Value = "AeroTorque_lbfFt_fPowerHorsepowerRevPerMin"
$ws.Range("D380").Value = "powerHorsepower"
$ws.Range("E380").Value = " revPerMin"
$ws.Range("H380").Value = 2
$ws.Range("I380").Formula = '="Torque in " & K380 & ", f(" & J380 & ")"'
$ws.Range("J380").Value = "horsepower, RPM"
$ws.Range("K380").Value = "lbf-ft"

$ws.Range("A381").Value = "TORQUE"
$ws.Range("B381").Value = "TORQUE"
$ws.Range("C381").Value = "AeroTorque_lbfFt_fPowerHorsepowerRevPerSec"
$ws.Range("D381").Value = "powerHorsepower"
$ws.Range("E381").Value = " revPerSec"
$ws.Range("H381").Value = 2
$ws.Range("I381").Formula = '="Torque in " & K381 & ", f(" & J381 & ")"'
$ws.Range("J381").Value = "horsepower, RPS"
$ws.Range("K381").Value = "lbf-ft"

$ws.Range("A382").Value = "TORQUE"
$ws.Range("B382").Value = "TORQUE"
$ws.Range("C382").Value = "AeroTorque_lbfInch_fPowerFtLbfPerSecRevPerMin"
$ws.Range("D382").Value = "powerFtLbfPerSec"
$ws.Range("E382").Value = " revPerMin"
$ws.Range("H382").Value = 2
$ws.Range("I382").Formula = '="Torque in " & K382 & ", f(" & J382 & ")"'
$ws.Range("J382").Value = "lbf-ft/sec, RPM"
$ws.Range("K382").Value = "lbf-inch"

$ws.Range("A383").Value = "TORQUE"
$ws.Range("B383").Value = "TORQUE"
$ws.Range("C383").Value = "AeroTorque_lbfInch_fPowerFtLbfPerSecRevPerSec"
$ws.Range("D383").Value = "powerFtLbfPerSec"
$ws.Range("E383").Value = " revPerSec"
$ws.Range("H383").Value = 2
$ws.Range("I383").Formula = '="Torque in " & K383 & ", f(" & J383 & ")"'
$ws.Range("J383").Value = "lbf-ft/sec, RPS"
$ws.Range("K383").Value = "lbf-inch"

$ws.Range("A384").Value = "TORQUE"
$ws.Range("B384").Value = "TORQUE"
$ws.Range("C384").Value = "AeroTorque_lbfInch_fPowerHorsepowerRevPerMin"
$ws.Range("D384").Value = "powerHorsepower"
$ws.Range("E384").Value = " revPerMin"
$ws.Range("H384").Value = 2
$ws.Range("I384").Formula = '="Torque in " & K384 & ", f(" & J384 & ")"'
$ws.Range("J384").Value = "horsepower, RPM"
$ws.Range("K384").Value = "lbf-inch"

$ws.Range("A385").Value = "TORQUE"
$ws.Range("B385").Value = "TORQUE"
$ws.Range("C385").Value = "AeroTorque_lbfInch_fPowerHorsepowerRevPerSec"
$ws.Range("D385").Value = "powerHorsepower"
$ws.Range("E385").Value = " revPerSec"
$ws.Range("H385").Value = 2
$ws.Range("I385").Formula = '="Torque in " & K385 & ", f(" & J385 & ")"'
$ws.Range("J385").Value = "horsepower, RPS"
$ws.Range("K385").Value = "lbf-inch"

$ws.Range("A386").Value = "TORQUE"
$ws.Range("B386").Value = "TORQUE"
$ws.Range("C386").Value = "AeroTorque_newtonMeter_fPowerWattsRevPerMin"
$ws.Range("D386").Value = "powerWatts"
$ws.Range("E386").Value = " revPerMin"
$ws.Range("H386").Value = 2
$ws.Range("I386").Formula = '="Torque in " & K386 & ", f(" & J386 & ")"'
$ws.Range("J386").Value = "watts, RPM"
$ws.Range("K386").Value = "newton-meter"

$ws.Range("A387").Value = "TORQUE"
$ws.Range("B387").Value = "TORQUE"
$ws.Range("C387").Value = "AeroTorque_newtonMeter_fPowerWattsRevPerSec"
$ws.Range("D387").Value = "powerWatts"
$ws.Range("E387").Value = " revPerSec"
$ws.Range("H387").Value = 2
$ws.Range("I387").Formula = '="Torque in " & K387 & ", f(" & J387 & ")"'
$ws.Range("J387").Value = "watts, RPS"
$ws.Range("K387").Value = "newton-meter"

$ws.Range("A388").Value = "TORQUE"
$ws.Range("B388").Value = "TORQUE"
$ws.Range("C388").Value = "AeroTorque_newtonMeter_fPowerKilowattsRevPerMin"
$ws.Range("D388").Value = "powerKilowatts"
$ws.Range("E388").Value = " revPerMin"
$ws.Range("H388").Value = 2
$ws.Range("I388").Formula = '="Torque in " & K388 & ", f(" & J388 & ")"'
$ws.Range("J388").Value = "kilowatts, RPM"
$ws.Range("K388").Value = "newton-meter"

$ws.Range("A389").Value = "TORQUE"
$ws.Range("B389").Value = "TORQUE"
$ws.Range("C389").Value = "AeroTorque_newtonMeter_fPowerKilowattsRevPerSec"
$ws.Range("D389").Value = "powerKilowatts"
$ws.Range("E389").Value = " revPerSec"
$ws.Range("H389").Value = 2
$ws.Range("I389").Formula = '="Torque in " & K389 & ", f(" & J389 & ")"'
$ws.Range("J389").Value = "kilowatts, RPS"
$ws.Range("K389").Value = "newton-meter"

# ---------------------------------------------------------------------
# 4) Column D needs to be a touch wider to fit "torqueNewtonMeter"/
#    "powerFtLbfPerSec" etc. without truncation.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 21

$ws.Calculate()